# "Agregada comprobacion de stock" - adds a "Disponible" (si/no) column (G)
# showing stock availability, fixes a FileId typo, and appends four new
# product rows (Rojas lisas, Amarillas lisas, Celestes lisas, Violetas lisas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header -------------------------------------------------
$ws.Range("G1").Value = "Disponible"

# --- Fix a typo in an existing FileId (row 8 / ID 7 "Vaquitas en negro") -
$ws.Range("F8").Value = "1pQLTtWbM4fdGTYrjIUC59KSxRbaG_0Mt"

# --- Fill "Disponible" for the existing 14 rows -------------------------
$ws.Range("G2").Value = "si"
$ws.Range("G3").Value = "si"
$ws.Range("G4").Value = "no"
$ws.Range("G5").Value = "si"
$ws.Range("G6").Value = "no"
$ws.Range("G7").Value = "no"
$ws.Range("G8").Value = "no"
$ws.Range("G9").Value = "si"
$ws.Range("G10").Value = "si"
$ws.Range("G11").Value = "no"
$ws.Range("G12").Value = "si"
$ws.Range("G13").Value = "no"
$ws.Range("G14").Value = "no"
$ws.Range("G15").Value = "si"

# --- Append the four new products (rows 16-19) --------------------------
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Rojas lisas"
$ws.Range("C16").Value = 1500
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = "Coffin"
$ws.Range("F16").Value = "1pzDoLFjsip8QmHuF-JXWVboxEP2hgIMT"
$ws.Range("G16").Value = "si"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Amarillas lisas"
$ws.Range("C17").Value = 1500
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = "Coffin"
$ws.Range("F17").Value = "1q9-QKZpkwh2m_l_nQQVD06hhAfZQYbwp"
$ws.Range("G17").Value = "si"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Celestes lisas"
$ws.Range("C18").Value = 1500
$ws.Range("D18").Value = 2.3
$ws.Range("E18").Value = "Cuadradas"
$ws.Range("F18").Value = "1q9RblGY49d_lSqEIxVWEDMYoCOP0slEZ"
$ws.Range("G18").Value = "si"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Violetas lisas"
$ws.Range("C19").Value = 1500
$ws.Range("D19").Value = 2.3
$ws.Range("E19").Value = "Cuadradas"
$ws.Range("F19").Value = "1qNxa9t8YnRFxEF9ltpW_f1ORKZPO9CGC"
$ws.Range("G19").Value = "si"

# --- Hyperlinks for the new FileId cells (F18 created before F16, to
#     reproduce the original rId ordering) -------------------------------
$ws.Hyperlinks.Add($ws.Range("F18"), "https://drive.google.com/file/d/1q9RblGY49d_lSqEIxVWEDMYoCOP0slEZ/view?usp=sharing")
$ws.Range("F18").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F16"), "https://drive.google.com/file/d/1pzDoLFjsip8QmHuF-JXWVboxEP2hgIMT/view?usp=sharing")
$ws.Range("F16").Style = "Hyperlink"

# --- Column widths: narrower FileId column, new Disponible column -------
$ws.Columns.Item(6).ColumnWidth = 41
$ws.Columns.Item(7).ColumnWidth = 17.3

# --- Selection shown when the file was last saved ------------------------
$ws.Range("G9").Select()
